$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configs")

# 1. Insert a new row above row 3; this shifts existing rows/merges down by one.
$ws.Rows("3:3").Insert()

# 2. Populate the new row 3 with the new header values.
$ws.Range("B3").Value = "properties"
$ws.Range("C3").Value = "origin"
$ws.Range("D3").Value = "Deviation"

Write-Host "stage1 done"
